# Updated cryptos list on Wed Nov 22 11:47:00 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.573.79"
$ws.Range("E2").Value = "  -1.79%  "

$ws.Range("D3").Value = "2.022.17"
$ws.Range("E3").Value = "  +0.86%  "

$ws.Range("D5").Value = "234.78"
$ws.Range("E5").Value = "  -9.34%  "

$ws.Range("E6").Value = "  -2.63%  "

$ws.Range("E7").Value = "  +0.01%  "

$ws.Range("D8").Value = "54.89"
$ws.Range("E8").Value = "  -2.83%  "

$ws.Range("E9").Value = "  -2.76%  "

$ws.Range("D10").Value = "57.50"
$ws.Range("E10").Value = "  +2.58%  "

$ws.Range("E11").Value = "  -2.92%  "

$ws.Range("E12").Value = "  -0.56%  "

$ws.Range("D13").Value = "2.320.16"
$ws.Range("E13").Value = "  +0.75%  "

$ws.Range("D14").Value = "14.20"
$ws.Range("E14").Value = "  -0.05%  "

$ws.Range("D15").Value = "20.14"
$ws.Range("E15").Value = "  -7.04%  "

$ws.Range("D16").Value = "0.764"
$ws.Range("E16").Value = "  -3.59%  "

$ws.Range("E17").Value = "  -2.01%  "

$ws.Range("D18").Value = "2.022.86"
$ws.Range("E18").Value = "  +0.65%  "

$ws.Range("D19").Value = "36.477.37"
$ws.Range("E19").Value = "  -2.15%  "

$ws.Range("D20").Value = "67.70"
$ws.Range("E20").Value = "  -3.38%  "

$ws.Range("D21").Value = "0.0₃0796"
$ws.Range("E21").Value = "  -4.34%  "

$ws.Range("D22").Value = "5.40"
$ws.Range("E22").Value = "  +6.02%  "

$ws.Range("D23").Value = "220.36"
$ws.Range("E23").Value = "  -5.65%  "

$ws.Range("E24").Value = "  -0.01%  "

$ws.Range("E25").Value = "  +1.36%  "

$ws.Range("D26").Value = "2.40"
$ws.Range("E26").Value = "  -6.54%  "

$ws.Range("D27").Value = "163.14"
$ws.Range("E27").Value = "  -1.05%  "

$ws.Range("D28").Value = "8.59"
$ws.Range("E28").Value = "  -3.87%  "

# Row 29/30 swap: ImmutableX <-> Kaspa
$ws.Range("B29").Value = "Kaspa"
$ws.Range("C29").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D29").Value = "0.129"
$ws.Range("E29").Value = "  +0.24%  "

$ws.Range("B30").Value = "ImmutableX"
$ws.Range("C30").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D30").Value = "1.38"
$ws.Range("E30").Value = "  +4.94%  "

$ws.Range("D31").Value = "18.94"
$ws.Range("E31").Value = "  -2.99%  "

$ws.Range("E32").Value = "  -1.84%  "

$ws.Range("E33").Value = "  -4.76%  "

$ws.Range("E34").Value = "  -5.29%  "

$ws.Range("E35").Value = "  +4.21%  "

$ws.Range("D36").Value = "4.24"
$ws.Range("E36").Value = "  -4.26%  "

$ws.Range("E37").Value = "  -0.11%  "

$ws.Range("D38").Value = "1.77"
$ws.Range("E38").Value = "  -2.30%  "

$ws.Range("E39").Value = "  -3.18%  "

$ws.Range("D40").Value = "5.68"
$ws.Range("E40").Value = "  +4.45%  "

$ws.Range("D41").Value = "2.99"
$ws.Range("E41").Value = "  -2.13%  "

# Row 42/43 swap: Maker <-> Cronos
$ws.Range("B42").Value = "Cronos"
$ws.Range("C42").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D42").Value = "0.0949"
$ws.Range("E42").Value = "  +2.78%  "

$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").Value = "1.455.36"
$ws.Range("E43").Value = "  +1.38%  "

$ws.Range("E45").Value = "  -3.27%  "

# Row 46/47 swap: TrustWalletToken <-> Aave
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").Value = "90.14"
$ws.Range("E46").Value = "  +1.13%  "

$ws.Range("B47").Value = "TrustWalletToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D47").Value = "1.10"
$ws.Range("E47").Value = "  -6.41%  "

$ws.Range("E48").Value = "  -1.19%  "

$ws.Range("E49").Value = "  -1.02%  "

$ws.Range("E50").Value = "  -1.79%  "

$ws.Range("D51").Value = "6.86"
$ws.Range("E51").Value = "  -1.43%  "
